# router revised and refined
#
# Applies the "routes" worksheet edits:
#  - drop the shared string "erro" (its two usages, R6/T6, become blank
#    green cells) and introduce a new shared string "desistir" (placed in S9)
#  - recolor / repurpose a handful of helper cells in columns Q-T of the
#    routes sheet (they track "a editar" / "remover link" markers)
#  - retarget R3 from "s/ conteudo" to "remover link"
#  - move the active selection to Q3
#
# Colors used by the sheet (Interior.Color, BGR long values):
#   GREEN  (92D050) = 5296274
#   ORANGE (FFC000)  = 49407
#   RED    (FF0000)  = 255

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1 -----------------------------------------------------------
# S1 picks up the red "needs work" marker, T1 moves from blue to green.
$ws.Range("S1").Interior.Color = 255
$ws.Range("T1").Interior.Color = 5296274

# --- Row 2 -----------------------------------------------------------
# Q2 / R2 / T2 become "a editar" markers (green); S2 stays as-is.
$ws.Range("Q2").Value = "a editar"
$ws.Range("Q2").Interior.Color = 5296274
$ws.Range("R2").Value = "a editar"
$ws.Range("R2").Interior.Color = 5296274
$ws.Range("T2").Value = "a editar"
$ws.Range("T2").Interior.Color = 5296274

# --- Row 3 -----------------------------------------------------------
# R3 used to say "s/ conteudo" (orange) -- now it is a red "remover link".
$ws.Range("R3").Value = "remover link"
$ws.Range("R3").Interior.Color = 255

# --- Row 4 -----------------------------------------------------------
# New red "remover link" / blank markers in R4 and T4.
$ws.Range("R4").Value = "remover link"
$ws.Range("R4").Interior.Color = 255
$ws.Range("T4").Interior.Color = 255

# --- Row 5 -----------------------------------------------------------
# New green "a editar" markers in R5 and T5.
$ws.Range("R5").Value = "a editar"
$ws.Range("R5").Interior.Color = 5296274
$ws.Range("T5").Value = "a editar"
$ws.Range("T5").Interior.Color = 5296274

# --- Row 6 -----------------------------------------------------------
# R6 / T6 carried the now-retired "erro" label (orange) -- clear them and
# recolor green so the "erro" shared string drops out entirely.
$ws.Range("R6").ClearContents()
$ws.Range("T6").ClearContents()
$ws.Range("R6").Interior.Color = 5296274
$ws.Range("T6").Interior.Color = 5296274

# --- Row 7 -----------------------------------------------------------
# New blank green markers in R7 and T7.
$ws.Range("R7").Interior.Color = 5296274
$ws.Range("T7").Interior.Color = 5296274

# --- Row 9 -----------------------------------------------------------
# Brand new "desistir" (give-up) marker, orange, next to the notes row.
$ws.Range("S9").Value = "desistir"
$ws.Range("S9").Interior.Color = 49407

# --- Row 13 ------------------------------------------------------------
# A trailing helper cell next to the W13 marker.
$ws.Range("P13").Value = ""

# --- Selection ---------------------------------------------------------
$ws.Range("Q3").Select()
